# Updated capital structure database
# Applies the diff described in the commit: updates numeric metrics for
# the Jamaica "Investments & Asset Management" rows (2-6), renames three
# companies (rows shifted identity), and clears cells that no longer have values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company name (text) updates ---
$ws.Range("B4").Value = "1834 Investments Limited (JMSE:1834)"
$ws.Range("B5").Value = "QWI Investments Limited (JMSE:QWI)"
$ws.Range("B6").Value = "SSL Venture Capital Jamaica Limited (JMSE:SSLVC)"

# --- Numeric cell updates, grouped by row ---
# Row 2
$ws.Range("G2").Value = 0.3790477636106051
$ws.Range("H2").Value = 0.3790477636106051
$ws.Range("I2").Value = 0.1777727180732645
$ws.Range("J2").Value = 0.1535198872220369
$ws.Range("K2").Value = 7.297
$ws.Range("L2").Value = 0.1846033191661607
$ws.Range("M2").Value = 8.751
$ws.Range("N2").Value = 0.04864909939959973
$ws.Range("O2").Value = 1.199259969850623
$ws.Range("P2").Value = 8.751
$ws.Range("Q2").Value = 0.04864909939959973
$ws.Range("R2").Value = 1.199259969850623
$ws.Range("U2").Value = 125.881
$ws.Range("V2").Value = 0.6998054258394485
$ws.Range("W2").Value = 0.7381380519462157
$ws.Range("X2").Value = 0.07054579870886214
$ws.Range("Y2").Value = 0.6675922532373536
$ws.Range("Z2").Value = 1.589129211224572
$ws.Range("AA2").Value = -0.1266461201988023
$ws.Range("AB2").Value = 0.06825741611347161
$ws.Range("AC2").Value = -0.193366653933624
$ws.Range("AD2").Value = 151.958
$ws.Range("AE2").Value = 0.0
$ws.Range("AF2").Value = 151.958
$ws.Range("AG2").Value = 26.077
$ws.Range("AH2").Value = 0.457928266202183
$ws.Range("AI2").Value = 0.5057377158299718
$ws.Range("AJ2").Value = 0.1266138077365664
$ws.Range("AK2").Value = 0.1493639274401874
$ws.Range("AL2").Value = 11.033
$ws.Range("AM2").Value = 11.014
$ws.Range("AN2").Value = 12.3623494956069
$ws.Range("AO2").Value = 0.6369074594398622
$ws.Range("AP2").Value = 2.121461112918972
$ws.Range("AQ2").Value = 0.6380061739604139

# Row 3
$ws.Range("G3").Value = 0.3729216152019002
$ws.Range("H3").Value = 0.3729216152019002
$ws.Range("I3").Value = 0.2731591448931116
$ws.Range("J3").Value = 0.247209026128266
$ws.Range("K3").Value = 11.0
$ws.Range("L3").Value = 0.2612826603325415
$ws.Range("M3").Value = 8.08
$ws.Range("N3").Value = 0.04969249692496925
$ws.Range("O3").Value = 0.7345454545454545
$ws.Range("P3").Value = 8.08
$ws.Range("Q3").Value = 0.04969249692496925
$ws.Range("R3").Value = 0.7345454545454545
$ws.Range("U3").Value = 125.1
$ws.Range("V3").Value = 0.7693726937269373
$ws.Range("X3").Value = 0.1028660256873606
$ws.Range("AB3").Value = 0.08176835819041203
$ws.Range("AD3").Value = 150.4
$ws.Range("AE3").Value = 0.0
$ws.Range("AF3").Value = 150.4
$ws.Range("AG3").Value = 25.30000000000001
$ws.Range("AH3").Value = 0.4805111821086262
$ws.Range("AI3").Value = 0.5375268048606148
$ws.Range("AJ3").Value = 0.1346460883448644
$ws.Range("AK3").Value = 0.1635423400129283
$ws.Range("AL3").Value = 10.8
$ws.Range("AM3").Value = 10.8
$ws.Range("AN3").Value = 11.48091603053435
$ws.Range("AO3").Value = 1.064814814814815
$ws.Range("AP3").Value = 1.931297709923665
$ws.Range("AQ3").Value = 1.064814814814815

# Row 4
$ws.Range("G4").Value = -1.610062893081761
$ws.Range("H4").Value = -1.610062893081761
$ws.Range("I4").Value = -2.327044025157233
$ws.Range("J4").Value = -1.27823545043848
$ws.Range("K4").Value = 0.078
$ws.Range("L4").Value = 0.4905660377358491
$ws.Range("M4").Value = 0.671
$ws.Range("N4").Value = 0.08192918192918194
$ws.Range("O4").Value = 8.602564102564104
$ws.Range("P4").Value = 0.671
$ws.Range("Q4").Value = 0.08192918192918194
$ws.Range("R4").Value = 8.602564102564104
$ws.Range("T4").Value = 0.0
$ws.Range("U4").Value = 0.75
$ws.Range("V4").Value = 0.09157509157509158
$ws.Range("W4").Value = 0.007289719626168225
$ws.Range("X4").Value = 0.06454268619014099
$ws.Range("Y4").Value = -0.05725296656397276
$ws.Range("Z4").Value = 0.01730141458106638
$ws.Range("AA4").Value = -0.02211528146025226
$ws.Range("AB4").Value = 0.06454268619014099
$ws.Range("AC4").Value = -0.08665796765039324
$ws.Range("AD4").Value = 0.0
$ws.Range("AF4").Value = 0.0
$ws.Range("AG4").Value = -0.75
$ws.Range("AH4").Value = 0.0
$ws.Range("AI4").Value = 0.0
$ws.Range("AJ4").Value = -0.1008064516129032
$ws.Range("AK4").Value = -0.07772020725388601
$ws.Range("AL4").Value = 0.0
$ws.Range("AM4").Value = -0.005
$ws.Range("AN4").Value = -0.0
$ws.Range("AP4").Value = 2.124645892351275
$ws.Range("AQ4").Value = 74.0

# Row 5
$ws.Range("G5").Value = -0.0
$ws.Range("H5").Value = -0.0
$ws.Range("I5").Value = 1.087087087087087
$ws.Range("J5").Value = 1.087087087087087
$ws.Range("K5").Value = -2.81
$ws.Range("L5").Value = 0.8438438438438438
$ws.Range("M5").Value = -0.0
$ws.Range("N5").Value = -0.0
$ws.Range("O5").Value = 0.0
$ws.Range("P5").Value = -0.0
$ws.Range("Q5").Value = -0.0
$ws.Range("R5").Value = 0.0
$ws.Range("U5").Value = 0.007
$ws.Range("V5").Value = 0.0009370816599732263
$ws.Range("W5").Value = -0.202158273381295
$ws.Range("X5").Value = 0.07225225785723899
$ws.Range("Y5").Value = -0.274410531238534
$ws.Range("Z5").Value = -0.2126572578070119
$ws.Range("AA5").Value = -0.2311769589373523
$ws.Range("AB5").Value = 0.06889838127950233
$ws.Range("AC5").Value = -0.3000753402168547
$ws.Range("AD5").Value = 1.39
$ws.Range("AF5").Value = 1.39
$ws.Range("AG5").Value = 1.383
$ws.Range("AH5").Value = 0.1568848758465011
$ws.Range("AI5").Value = 0.1169049621530698
$ws.Range("AJ5").Value = 0.156218231108099
$ws.Range("AK5").Value = 0.1163847513254229
$ws.Range("AL5").Value = 0.141
$ws.Range("AM5").Value = 0.127
$ws.Range("AO5").Value = -25.67375886524823
$ws.Range("AQ5").Value = -28.50393700787402

# Row 6
$ws.Range("G6").Value = -0.7696160267111853
$ws.Range("H6").Value = -0.7696160267111853
$ws.Range("I6").Value = -0.8063439065108514
$ws.Range("J6").Value = -0.8063439065108514
$ws.Range("K6").Value = -0.971
$ws.Range("L6").Value = -1.621035058430718
$ws.Range("M6").Value = -0.0
$ws.Range("N6").Value = -0.0
$ws.Range("O6").Value = 0.0
$ws.Range("P6").Value = -0.0
$ws.Range("Q6").Value = -0.0
$ws.Range("R6").Value = 0.0
$ws.Range("U6").Value = 0.024
$ws.Range("V6").Value = 0.01481481481481481
$ws.Range("W6").Value = 1.468986384266263
$ws.Range("X6").Value = 0.0688393395604853
$ws.Range("Y6").Value = 1.400147044705778
$ws.Range("Z6").Value = 23.95999999999998
$ws.Range("AA6").Value = -19.31999999999998
$ws.Range("AB6").Value = 0.0676164509474409
$ws.Range("AC6").Value = -19.38761645094742
$ws.Range("AD6").Value = 0.168
$ws.Range("AF6").Value = 0.168
$ws.Range("AG6").Value = 0.144
$ws.Range("AH6").Value = 0.09395973154362416
$ws.Range("AI6").Value = -0.1035758323057953
$ws.Range("AJ6").Value = 0.08163265306122448
$ws.Range("AK6").Value = -0.08748481166464157
$ws.Range("AL6").Value = 0.092
$ws.Range("AM6").Value = 0.092
$ws.Range("AN6").Value = -0.3692307692307693
$ws.Range("AO6").Value = -5.25
$ws.Range("AP6").Value = -0.3164835164835165
$ws.Range("AQ6").Value = -5.25

# --- Clear cells that are no longer populated in the updated dataset ---
$cellsToClear = @("D2", "E2", "Z3", "AA3", "AC3", "D4", "AO4", "D5", "E5", "T5", "T6")
foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}
